$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (column D) values, forcing text storage to avoid numeric auto-conversion
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "24.616.23"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.693.79"
$ws.Range("D3").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.58"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3943"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4015"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.519"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.001"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "52.31"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08756"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.217"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.24"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.237"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001311"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.698.77"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "99.82"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.07067"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.66"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.044"
$ws.Range("D21").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.18"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "24.621.54"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.125"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.341"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.74"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "161.88"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "136.87"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.203"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.448"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.884.46"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.077"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08597"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "7.177"
$ws.Range("D35").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2732"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.924"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "14.37"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.09116"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.02730"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.480"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.7621"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.597"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.7151"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "15.63"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.222"
$ws.Range("D47").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "140.91"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.323"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07982"
$ws.Range("D51").Style = "Normal"

# Update Volume(1h) (column E) values
$ws.Range("E2").Value = "  +2.30%  "
$ws.Range("E3").Value = "  +1.45%  "
$ws.Range("E4").Value = "  +0.62%  "
$ws.Range("E5").Value = "  +1.81%  "
$ws.Range("E6").Value = "  +0.53%  "
$ws.Range("E7").Value = "  +1.41%  "
$ws.Range("E8").Value = "  +0.60%  "
$ws.Range("E9").Value = "  +3.52%  "
$ws.Range("E10").Value = "  +0.65%  "
$ws.Range("E11").Value = "  -2.73%  "
$ws.Range("E12").Value = "  +0.80%  "
$ws.Range("E13").Value = "  +6.01%  "
$ws.Range("E15").Value = "  +12.55%  "
$ws.Range("E16").Value = "  -0.20%  "
$ws.Range("E17").Value = "  +1.69%  "
$ws.Range("E18").Value = "  +0.52%  "
$ws.Range("E19").Value = "  +2.79%  "
$ws.Range("E20").Value = "  +2.39%  "
$ws.Range("E21").Value = "  +5.94%  "
$ws.Range("E22").Value = "  +0.23%  "
$ws.Range("E23").Value = "  +2.52%  "
$ws.Range("E24").Value = "  +2.37%  "
$ws.Range("E25").Value = "  +8.55%  "
$ws.Range("E26").Value = "  +1.71%  "
$ws.Range("E27").Value = "  +4.48%  "
$ws.Range("E28").Value = "  +2.02%  "
$ws.Range("E29").Value = "  +5.05%  "
$ws.Range("E30").Value = "  +1.97%  "
$ws.Range("E31").Value = "  +4.93%  "
$ws.Range("E32").Value = "  +1.60%  "
$ws.Range("E33").Value = "  -4.67%  "
$ws.Range("E34").Value = "  +0.09%  "
$ws.Range("E35").Value = "  +6.95%  "
$ws.Range("E36").Value = "  +9.69%  "
$ws.Range("E37").Value = "  +2.71%  "
$ws.Range("E38").Value = "  -0.62%  "
$ws.Range("E39").Value = "  -1.51%  "
$ws.Range("E40").Value = "  +3.38%  "
$ws.Range("E41").Value = "  +7.71%  "
$ws.Range("E42").Value = "  +1.21%  "
$ws.Range("E43").Value = "  +0.13%  "
$ws.Range("E44").Value = "  +7.01%  "
$ws.Range("E45").Value = "  +0.55%  "
$ws.Range("E46").Value = "  +3.89%  "
$ws.Range("E47").Value = "  +2.67%  "
$ws.Range("E48").Value = "  +0.35%  "
$ws.Range("E49").Value = "  +0.27%  "
$ws.Range("E50").Value = "  +7.92%  "
$ws.Range("E51").Value = "  +1.95%  "
